# "Generate Report for Handback" - update the handback status report
# after a new handback run:
#   * the zh-cn / de-de locales are no longer in sync with en-US, so the
#     status text changes everywhere it is shown (Overview + per-locale
#     sheets)
#   * the b5aba85f... file just got handed back again in both locales, so
#     its "Correspond Handback DateTime" is refreshed on each locale sheet
#   * the Status columns got a bit wider to fit the new (longer) text

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Handed back: not in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns for both files ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn detail sheet: Status column ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# --- de-de detail sheet: Status column ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Refresh the Correspond Handback DateTime for the b5aba85f... file ---
# (new handback pass timestamps for each locale)
$wsZhCn.Range("K3").Value = "2016-11-14 07:11:59"
$wsDeDe.Range("K3").Value = "2016-11-14 07:12:17"

# --- Widen the Status / locale columns so the longer text still fits ---
# (target character width ~33.46; the COM width setter here snaps to a
# pixel grid, so feed it the matching pre-snap value)
$newColumnWidth = 32.62688700358076

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
